$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 7.157807
$ws.Range("H2").Value = 14.315614
$ws.Range("I2").Value = 0.21112489480340316
$ws.Range("J2").Value = 0.16565465444792946
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.569834
$ws.Range("N2").Value = 1.139668
$ws.Range("O2").Value = 0.0076706049295135375
$ws.Range("P2").Value = 0.005687058042690201
$ws.Range("Q2").Value = 4.078761794038
$ws.Range("R2").Value = 16.315047176152
$ws.Range("S2").Value = 0.0016194556588220114
$ws.Range("T2").Value = 0.0009420876348871633

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 7.157807
$ws.Range("H3").Value = 14.315614
$ws.Range("I3").Value = 0.21112489480340316
$ws.Range("J3").Value = 0.16565465444792946
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.102415
$ws.Range("N3").Value = 150.307245
$ws.Range("O3").Value = 0.674434715161842
$ws.Range("P3").Value = 0.7500482829664924
$ws.Range("Q3").Value = 358.623416803905
$ws.Range("R3").Value = 2151.74050082343
$ws.Range("S3").Value = 0.14238995829030704
$ws.Range("T3").Value = 0.12424898913407711

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 7.157807
$ws.Range("H4").Value = 14.315614
$ws.Range("I4").Value = 0.21112489480340316
$ws.Range("J4").Value = 0.16565465444792946
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4051446666666667
$ws.Range("N4").Value = 1.2154340000000001
$ws.Range("O4").Value = 0.005453701739979455
$ws.Range("P4").Value = 0.00606513800954236
$ws.Range("Q4").Value = 2.899947331079334
$ws.Range("R4").Value = 17.399683986476003
$ws.Range("S4").Value = 0.0011514122061422993
$ws.Range("T4").Value = 0.0010047183411497423

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 7.157807
$ws.Range("H5").Value = 14.315614
$ws.Range("I5").Value = 0.21112489480340316
$ws.Range("J5").Value = 0.16565465444792946
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2841426666666667
$ws.Range("N5").Value = 0.8524280000000002
$ws.Range("O5").Value = 0.0038248790693753898
$ws.Range("P5").Value = 0.004253701528176911
$ws.Range("Q5").Value = 2.0338383684653336
$ws.Range("R5").Value = 12.203030210792003
$ws.Range("S5").Value = 0.0008075271911576178
$ws.Range("T5").Value = 0.0007046454567747757

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 7.157807
$ws.Range("H6").Value = 14.315614
$ws.Range("I6").Value = 0.21112489480340316
$ws.Range("J6").Value = 0.16565465444792946
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.029030333333333
$ws.Range("N6").Value = 3.087090999999999
$ws.Range("O6").Value = 0.013851902742703355
$ws.Range("P6").Value = 0.01540489484662773
$ws.Range("Q6").Value = 7.365600523145664
$ws.Range("R6").Value = 44.19360313887399
$ws.Range("S6").Value = 0.0029244815093802173
$ws.Range("T6").Value = 0.002551892532624806

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 7.157807
$ws.Range("H7").Value = 14.315614
$ws.Range("I7").Value = 0.21112489480340316
$ws.Range("J7").Value = 0.16565465444792946
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 21.8974465
$ws.Range("N7").Value = 43.794893
$ws.Range("O7").Value = 0.29476419635658624
$ws.Range("P7").Value = 0.21854092460647032
$ws.Range("Q7").Value = 156.7376958398255
$ws.Range("R7").Value = 626.950783359302
$ws.Range("S7").Value = 0.06223205994759394
$ws.Range("T7").Value = 0.03620232134841584

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.722124333333333
$ws.Range("H8").Value = 11.166372999999998
$ws.Range("I8").Value = 0.10978685347623748
$ws.Range("J8").Value = 0.12921287628680747
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 0.569834
$ws.Range("N8").Value = 1.139668
$ws.Range("O8").Value = 0.0076706049295135375
$ws.Range("P8").Value = 0.005687058042690201
$ws.Range("Q8").Value = 2.1209929973606663
$ws.Range("R8").Value = 12.725957984163998
$ws.Range("S8").Value = 0.0008421315794706077
$ws.Range("T8").Value = 0.0007348411273060225

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.722124333333333
$ws.Range("H9").Value = 11.166372999999998
$ws.Range("I9").Value = 0.10978685347623748
$ws.Range("J9").Value = 0.12921287628680747
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 50.102415
$ws.Range("N9").Value = 150.307245
$ws.Range("O9").Value = 0.674434715161842
$ws.Range("P9").Value = 0.7500482829664924
$ws.Range("Q9").Value = 186.487418030265
$ws.Range("R9").Value = 1678.3867622723847
$ws.Range("S9").Value = 0.07404406525276111
$ws.Range("T9").Value = 0.09691589599608175

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.722124333333333
$ws.Range("H10").Value = 11.166372999999998
$ws.Range("I10").Value = 0.10978685347623748
$ws.Range("J10").Value = 0.12921287628680747
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.4051446666666667
$ws.Range("N10").Value = 1.2154340000000001
$ws.Range("O10").Value = 0.005453701739979455
$ws.Range("P10").Value = 0.00606513800954236
$ws.Range("Q10").Value = 1.5079988223202223
$ws.Range("R10").Value = 13.571989400882
$ws.Range("S10").Value = 0.0005987447538302259
$ws.Range("T10").Value = 0.0007836939272894107

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.722124333333333
$ws.Range("H11").Value = 11.166372999999998
$ws.Range("I11").Value = 0.10978685347623748
$ws.Range("J11").Value = 0.12921287628680747
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.2841426666666667
$ws.Range("N11").Value = 0.8524280000000002
$ws.Range("O11").Value = 0.0038248790693753898
$ws.Range("P11").Value = 0.004253701528176911
$ws.Range("Q11").Value = 1.0576143337382222
$ws.Range("R11").Value = 9.518529003644
$ws.Range("S11").Value = 0.0004199214379538435
$ws.Range("T11").Value = 0.0005496330093213271

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.722124333333333
$ws.Range("H12").Value = 11.166372999999998
$ws.Range("I12").Value = 0.10978685347623748
$ws.Range("J12").Value = 0.12921287628680747
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.029030333333333
$ws.Range("N12").Value = 3.087090999999999
$ws.Range("O12").Value = 0.013851902742703355
$ws.Range("P12").Value = 0.01540489484662773
$ws.Range("Q12").Value = 3.8301788434381097
$ws.Range("R12").Value = 34.47160959094298
$ws.Range("S12").Value = 0.0015207568167802652
$ws.Range("T12").Value = 0.001990510772028587

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.722124333333333
$ws.Range("H13").Value = 11.166372999999998
$ws.Range("I13").Value = 0.10978685347623748
$ws.Range("J13").Value = 0.12921287628680747
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 21.8974465
$ws.Range("N13").Value = 43.794893
$ws.Range("O13").Value = 0.29476419635658624
$ws.Range("P13").Value = 0.21854092460647032
$ws.Range("Q13").Value = 81.50501845551483
$ws.Range("R13").Value = 489.030110733089
$ws.Range("S13").Value = 0.032361233635441425
$ws.Range("T13").Value = 0.028238301454780367

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.570476333333333
$ws.Range("H14").Value = 4.711428999999999
$ws.Range("I14").Value = 0.04632237927988757
$ws.Range("J14").Value = 0.05451880324175782
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 0.569834
$ws.Range("N14").Value = 1.139668
$ws.Range("O14").Value = 0.0076706049295135375
$ws.Range("P14").Value = 0.005687058042690201
$ws.Range("Q14").Value = 0.8949108109286663
$ws.Range("R14").Value = 5.369464865571999
$ws.Range("S14").Value = 0.00035532067085110135
$ws.Range("T14").Value = 0.00031005159845388344

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.570476333333333
$ws.Range("H15").Value = 4.711428999999999
$ws.Range("I15").Value = 0.04632237927988757
$ws.Range("J15").Value = 0.05451880324175782
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.102415
$ws.Range("N15").Value = 150.307245
$ws.Range("O15").Value = 0.674434715161842
$ws.Range("P15").Value = 0.7500482829664924
$ws.Range("Q15").Value = 78.68465700034498
$ws.Range("R15").Value = 708.1619130031048
$ws.Range("S15").Value = 0.031241420675249782
$ws.Range("T15").Value = 0.04089173476086849

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.570476333333333
$ws.Range("H16").Value = 4.711428999999999
$ws.Range("I16").Value = 0.04632237927988757
$ws.Range("J16").Value = 0.05451880324175782
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.4051446666666667
$ws.Range("N16").Value = 1.2154340000000001
$ws.Range("O16").Value = 0.005453701739979455
$ws.Range("P16").Value = 0.00606513800954236
$ws.Range("Q16").Value = 0.6362701105762221
$ws.Range("R16").Value = 5.726430995185999
$ws.Range("S16").Value = 0.0002526284404787111
$ws.Range("T16").Value = 0.0003306640657763466

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.570476333333333
$ws.Range("H17").Value = 4.711428999999999
$ws.Range("I17").Value = 0.04632237927988757
$ws.Range("J17").Value = 0.05451880324175782
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.2841426666666667
$ws.Range("N17").Value = 0.8524280000000002
$ws.Range("O17").Value = 0.0038248790693753898
$ws.Range("P17").Value = 0.004253701528176911
$ws.Range("Q17").Value = 0.4462393332902222
$ws.Range("R17").Value = 4.016153999612
$ws.Range("S17").Value = 0.0001771774989513102
$ws.Range("T17").Value = 0.0002319067166638416

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1.570476333333333
$ws.Range("H18").Value = 4.711428999999999
$ws.Range("I18").Value = 0.04632237927988757
$ws.Range("J18").Value = 0.05451880324175782
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.029030333333333
$ws.Range("N18").Value = 3.087090999999999
$ws.Range("O18").Value = 0.013851902742703355
$ws.Range("P18").Value = 0.01540489484662773
$ws.Range("Q18").Value = 1.6160677847821103
$ws.Range("R18").Value = 14.544610063038993
$ws.Range("S18").Value = 0.0006416530925956197
$ws.Range("T18").Value = 0.0008398564311032662

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1.570476333333333
$ws.Range("H19").Value = 4.711428999999999
$ws.Range("I19").Value = 0.04632237927988757
$ws.Range("J19").Value = 0.05451880324175782
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 21.8974465
$ws.Range("N19").Value = 43.794893
$ws.Range("O19").Value = 0.29476419635658624
$ws.Range("P19").Value = 0.21854092460647032
$ws.Range("Q19").Value = 34.38942148868283
$ws.Range("R19").Value = 206.33652893209697
$ws.Range("S19").Value = 0.01365417890176104
$ws.Range("T19").Value = 0.011914589668891985

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 2.057700666666667
$ws.Range("H20").Value = 6.173102
$ws.Range("I20").Value = 0.06069342702127797
$ws.Range("J20").Value = 0.07143270827795596
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 0.569834
$ws.Range("N20").Value = 1.139668
$ws.Range("O20").Value = 0.0076706049295135375
$ws.Range("P20").Value = 0.005687058042690201
$ws.Range("Q20").Value = 1.1725478016893334
$ws.Range("R20").Value = 7.035286810135999
$ws.Range("S20").Value = 0.0004655553004984849
$ws.Range("T20").Value = 0.00040624195812329235

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 2.057700666666667
$ws.Range("H21").Value = 6.173102
$ws.Range("I21").Value = 0.06069342702127797
$ws.Range("J21").Value = 0.07143270827795596
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 50.102415
$ws.Range("N21").Value = 150.307245
$ws.Range("O21").Value = 0.674434715161842
$ws.Range("P21").Value = 0.7500482829664924
$ws.Range("Q21").Value = 103.09577274711
$ws.Range("R21").Value = 927.8619547239899
$ws.Range("S21").Value = 0.040933754165291644
$ws.Range("T21").Value = 0.053577980191527215

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 2.057700666666667
$ws.Range("H22").Value = 6.173102
$ws.Range("I22").Value = 0.06069342702127797
$ws.Range("J22").Value = 0.07143270827795596
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.4051446666666667
$ws.Range("N22").Value = 1.2154340000000001
$ws.Range("O22").Value = 0.005453701739979455
$ws.Range("P22").Value = 0.00606513800954236
$ws.Range("Q22").Value = 0.8336664506964446
$ws.Range("R22").Value = 7.502998056268001
$ws.Range("S22").Value = 0.00033100384855125976
$ws.Range("T22").Value = 0.0004332492341011819

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 2.057700666666667
$ws.Range("H23").Value = 6.173102
$ws.Range("I23").Value = 0.06069342702127797
$ws.Range("J23").Value = 0.07143270827795596
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.2841426666666667
$ws.Range("N23").Value = 0.8524280000000002
$ws.Range("O23").Value = 0.0038248790693753898
$ws.Range("P23").Value = 0.004253701528176911
$ws.Range("Q23").Value = 0.5846805546284446
$ws.Range("R23").Value = 5.262124991656001
$ws.Range("S23").Value = 0.0002321450186623488
$ws.Range("T23").Value = 0.0003038534203637568

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 2.057700666666667
$ws.Range("H24").Value = 6.173102
$ws.Range("I24").Value = 0.06069342702127797
$ws.Range("J24").Value = 0.07143270827795596
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 1.029030333333333
$ws.Range("N24").Value = 3.087090999999999
$ws.Range("O24").Value = 0.013851902742703355
$ws.Range("P24").Value = 0.01540489484662773
$ws.Range("Q24").Value = 2.1174364029202217
$ws.Range("R24").Value = 19.056927626281993
$ws.Range("S24").Value = 0.0008407194482201062
$ws.Range("T24").Value = 0.0011004133596317458

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 2.057700666666667
$ws.Range("H25").Value = 6.173102
$ws.Range("I25").Value = 0.06069342702127797
$ws.Range("J25").Value = 0.07143270827795596
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 21.8974465
$ws.Range("N25").Value = 43.794893
$ws.Range("O25").Value = 0.29476419635658624
$ws.Range("P25").Value = 0.21854092460647032
$ws.Range("Q25").Value = 45.058390261347675
$ws.Range("R25").Value = 270.35034156808604
$ws.Range("S25").Value = 0.017890249240054116
$ws.Range("T25").Value = 0.015610970114208762

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 11.26174566666667
$ws.Range("H26").Value = 33.78523700000001
$ws.Range("I26").Value = 0.33217364888124
$ws.Range("J26").Value = 0.39094947381763734
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 0.569834
$ws.Range("N26").Value = 1.139668
$ws.Range("O26").Value = 0.0076706049295135375
$ws.Range("P26").Value = 0.005687058042690201
$ws.Range("Q26").Value = 6.417325580219335
$ws.Range("R26").Value = 38.50395348131601
$ws.Range("S26").Value = 0.0025479728285629387
$ws.Range("T26").Value = 0.002223352349360097

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 11.26174566666667
$ws.Range("H27").Value = 33.78523700000001
$ws.Range("I27").Value = 0.33217364888124
$ws.Range("J27").Value = 0.39094947381763734
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 50.102415
$ws.Range("N27").Value = 150.307245
$ws.Range("O27").Value = 0.674434715161842
$ws.Range("P27").Value = 0.7500482829664924
$ws.Range("Q27").Value = 564.2406550157851
$ws.Range("R27").Value = 5078.165895142066
$ws.Range("S27").Value = 0.2240294402674888
$ws.Range("T27").Value = 0.2932309815635726

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 11.26174566666667
$ws.Range("H28").Value = 33.78523700000001
$ws.Range("I28").Value = 0.33217364888124
$ws.Range("J28").Value = 0.39094947381763734
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 0.4051446666666667
$ws.Range("N28").Value = 1.2154340000000001
$ws.Range("O28").Value = 0.005453701739979455
$ws.Range("P28").Value = 0.00606513800954236
$ws.Range("Q28").Value = 4.562636194206446
$ws.Range("R28").Value = 41.063725747858015
$ws.Range("S28").Value = 0.0018115760068789432
$ws.Range("T28").Value = 0.002371162513461938

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 11.26174566666667
$ws.Range("H29").Value = 33.78523700000001
$ws.Range("I29").Value = 0.33217364888124
$ws.Range("J29").Value = 0.39094947381763734
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 0.2841426666666667
$ws.Range("N29").Value = 0.8524280000000002
$ws.Range("O29").Value = 0.0038248790693753898
$ws.Range("P29").Value = 0.004253701528176911
$ws.Range("Q29").Value = 3.1999424450484457
$ws.Range("R29").Value = 28.799482005436015
$ws.Range("S29").Value = 0.0012705240370039048
$ws.Range("T29").Value = 0.0016629823742180434

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 11.26174566666667
$ws.Range("H30").Value = 33.78523700000001
$ws.Range("I30").Value = 0.33217364888124
$ws.Range("J30").Value = 0.39094947381763734
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 1.029030333333333
$ws.Range("N30").Value = 3.087090999999999
$ws.Range("O30").Value = 0.013851902742703355
$ws.Range("P30").Value = 0.01540489484662773
$ws.Range("Q30").Value = 11.588677897285223
$ws.Range("R30").Value = 104.298101075567
$ws.Range("S30").Value = 0.004601237077991829
$ws.Range("T30").Value = 0.006022535534505144

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 11.26174566666667
$ws.Range("H31").Value = 33.78523700000001
$ws.Range("I31").Value = 0.33217364888124
$ws.Range("J31").Value = 0.39094947381763734
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 21.8974465
$ws.Range("N31").Value = 43.794893
$ws.Range("O31").Value = 0.29476419635658624
$ws.Range("P31").Value = 0.21854092460647032
$ws.Range("Q31").Value = 246.60347323244025
$ws.Range("R31").Value = 1479.6208393946415
$ws.Range("S31").Value = 0.09791289866331357
$ws.Range("T31").Value = 0.08543845948251952

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 8.133334
$ws.Range("H32").Value = 16.266668
$ws.Range("I32").Value = 0.23989879653795387
$ws.Range("J32").Value = 0.18823148392791197
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 0.569834
$ws.Range("N32").Value = 1.139668
$ws.Range("O32").Value = 0.0076706049295135375
$ws.Range("P32").Value = 0.005687058042690201
$ws.Range("Q32").Value = 4.634650246555999
$ws.Range("R32").Value = 18.538600986223997
$ws.Range("S32").Value = 0.0018401688913083942
$ws.Range("T32").Value = 0.0010704833745597431

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 8.133334
$ws.Range("H33").Value = 16.266668
$ws.Range("I33").Value = 0.23989879653795387
$ws.Range("J33").Value = 0.18823148392791197
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 50.102415
$ws.Range("N33").Value = 150.307245
$ws.Range("O33").Value = 0.674434715161842
$ws.Range("P33").Value = 0.7500482829664924
$ws.Range("Q33").Value = 407.49967540161
$ws.Range("R33").Value = 2444.99805240966
$ws.Range("S33").Value = 0.1617960765107436
$ws.Range("T33").Value = 0.14118270132036528

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 8.133334
$ws.Range("H34").Value = 16.266668
$ws.Range("I34").Value = 0.23989879653795387
$ws.Range("J34").Value = 0.18823148392791197
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 0.4051446666666667
$ws.Range("N34").Value = 1.2154340000000001
$ws.Range("O34").Value = 0.005453701739979455
$ws.Range("P34").Value = 0.00606513800954236
$ws.Range("Q34").Value = 3.295176892318667
$ws.Range("R34").Value = 19.771061353912
$ws.Range("S34").Value = 0.0013083364840980165
$ws.Range("T34").Value = 0.0011416499277637408

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 8.133334
$ws.Range("H35").Value = 16.266668
$ws.Range("I35").Value = 0.23989879653795387
$ws.Range("J35").Value = 0.18823148392791197
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 0.2841426666666667
$ws.Range("N35").Value = 0.8524280000000002
$ws.Range("O35").Value = 0.0038248790693753898
$ws.Range("P35").Value = 0.004253701528176911
$ws.Range("Q35").Value = 2.311027211650667
$ws.Range("R35").Value = 13.866163269904002
$ws.Range("S35").Value = 0.000917583885646365
$ws.Range("T35").Value = 0.0008006805508351669

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 8.133334
$ws.Range("H36").Value = 16.266668
$ws.Range("I36").Value = 0.23989879653795387
$ws.Range("J36").Value = 0.18823148392791197
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 1.029030333333333
$ws.Range("N36").Value = 3.087090999999999
$ws.Range("O36").Value = 0.013851902742703355
$ws.Range("P36").Value = 0.01540489484662773
$ws.Range("Q36").Value = 8.369447397131331
$ws.Range("R36").Value = 50.21668438278798
$ws.Range("S36").Value = 0.0033230547977353172
$ws.Range("T36").Value = 0.0028996862167341817

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 8.133334
$ws.Range("H37").Value = 16.266668
$ws.Range("I37").Value = 0.23989879653795387
$ws.Range("J37").Value = 0.18823148392791197
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 21.8974465
$ws.Range("N37").Value = 43.794893
$ws.Range("O37").Value = 0.29476419635658624
$ws.Range("P37").Value = 0.21854092460647032
$ws.Range("Q37").Value = 178.099246131631
$ws.Range("R37").Value = 712.396984526524
$ws.Range("S37").Value = 0.07071357596842216
$ws.Range("T37").Value = 0.04113628253765384

